$d = $word.ActiveDocument

# Locate the last paragraph ("Sim pois ele se adapta...") which currently
# carries the _GoBack bookmark at its very end.
$old = "Sim pois ele se adapta melhor ao que o cliente pediu que é visualizar ao rápido um projeto já ser entregue rápido para ele conseguir já ir vendo com esta indo o projeto "
$lastPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq $old.TrimEnd()) {
        $lastPara = $cand
    }
}

# Remove the _GoBack bookmark from its current location; it will be
# re-created at the end of the newly appended content below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$newXml = "<w:p $ns/>" +
  "<w:p $ns><w:r><w:t>Como você organizaria a equipe de projeto?</w:t></w:r></w:p>" +
  "<w:p $ns>" +
    "<w:r><w:t>Teria um gestor de projetos(Scrum Master)que iria gerenciar o projeto e falar com o cliente para saber quais a exigências dele o que ele quer para o projeto definir o tempo para fazer cada coisa e também uma equipe de desenvolvimento(</w:t></w:r>" +
    "<w:proofErr $ns w:type='spellStart'/>" +
    "<w:r><w:t>Dev</w:t></w:r>" +
    "<w:proofErr $ns w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> Team)que com um grupo de pessoas boas em cada área para fazer com que o projeto saia no tempo em que foi acertado com o cliente</w:t></w:r>" +
    "<w:bookmarkStart $ns w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd $ns w:id='0'/>" +
  "</w:p>" +
  "<w:p $ns/>"

$endPos = $lastPara.Range.End
$insertAt = $d.Range($endPos, $endPos)
$insertAt.InsertXML($newXml)
